$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "/views/pages: This has the html code for the w|ebsite ..."
# used to be split into two runs with a leftover _GoBack bookmark sitting
# between the "w" and "ebsite" halves. Clean that up into one run holding
# the full, correctly spelled sentence.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$find = $d.Content.Find
$find.ClearFormatting()
[void]$find.Execute(
    "/views/pages: This has the html code for the website home page, Login page, and new user page.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/views/pages: This has the html code for the website home page, Login page, and new user page.",
    2
)

# ---------------------------------------------------------------------
# Edit 2: append two new notes at the very end of the document (still
# ahead of the final section break) -- a NOTE about the local map server
# and a TEST note about manual testing. Word's "last place you typed"
# bookmark (_GoBack) ends up on the new final paragraph.
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$xmlFragment = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p/>
<w:p>
  <w:r><w:t xml:space="preserve">NOTE: Since we run the map through a local server which KYLE got it may or may not show up when you try to run the </w:t></w:r>
  <w:r><w:t xml:space="preserve">website. </w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t>TEST: to test if it works, after creating the database  just m</w:t></w:r>
  <w:r><w:t xml:space="preserve">ake a new user and create tags. We did manual testing here. </w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$endRange.InsertXML($xmlFragment)
